$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "47 uF Ceramic Cap" row (row 9): new Digikey part num, quantities, unit cost, and a "NEW" note
$ws.Range("C9").Value = "587-1780-1-ND"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 0.7
$ws.Range("H9").Value = "NEW"

# Update selected cell to D7 (matches the saved selection in the sheet view)
$ws.Range("D7").Select()
